$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph (2nd paragraph in the body).
#    It sits right after the H1 title paragraph at the top of the doc.
# ---------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# ---------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play Abby and The Witch Slot for Free
#    - Review" right before the final "Prompt: ..." paragraph.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$lastRange = $lastPara.Range
$lastRange.Collapse(1)
$lastRange.InsertParagraphBefore()

$newPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$newRange = $newPara.Range
$newXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Abby and The Witch Slot for Free - Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newRange.InsertXML($newXml)

# ---------------------------------------------------------------------
# 3) Replace the body of the (now last) "Prompt: ..." paragraph with the
#    meta-description text (keeps its italic run formatting).
# ---------------------------------------------------------------------
$oldPrompt = "Prompt: Create a feature image for Abby & The Witch Design a cartoon-style feature image that includes a happy-looking Maya warrior wearing glasses. The image should also incorporate elements from the game " + [char]34 + "Abby & The Witch," + [char]34 + " such as Abby herself, the colorless world, and Baba Yaga's house and cemetery. Use bright colors to contrast the black and white world of the game and make the Maya warrior stand out. Feel free to add other magical elements to the image, like spells, potions, or magical creatures, to give it a more whimsical feel. The image should be eye-catching and convey the spirit of adventure and magic that the game offers to players."
$newDescription = "Abby and The Witch is a traditional online slot game with engaging graphics and storyline. Play for free and enjoy the free spin mode with a respectable RTP."

$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newDescription, 2)
